$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H2").Value = 3582
$ws.Range("I2").Value = 3927.75
$ws.Range("J2").Value = 2660
$ws.Range("K2").Value = 3927.75
$ws.Range("L2").Value = 2660
$ws.Range("M2").Value = -3814.75
$ws.Range("N2").Value = -2886
$ws.Range("H32").Value = 4947.406
$ws.Range("I32").Value = 3777.0566
$ws.Range("J32").Value = 8824.1875
$ws.Range("K32").Value = 3777.0566
$ws.Range("L32").Value = 8824.1875
$ws.Range("M32").Value = -3490.0566
$ws.Range("N32").Value = -9398.1875
$ws.Range("H61").Value = 5221.731
$ws.Range("I61").Value = 6036.476
$ws.Range("J61").Value = 1799.8
$ws.Range("K61").Value = 6036.476
$ws.Range("L61").Value = 1799.8
$ws.Range("M61").Value = -5824.476
$ws.Range("N61").Value = -2223.8
$ws.Range("H74").Value = 1426.1628
$ws.Range("I74").Value = 1302.9697
$ws.Range("K74").Value = 1302.9697
$ws.Range("M74").Value = -428.9697000000001
$ws.Range("H77").Value = 1426.1628
$ws.Range("I77").Value = 1302.9697
$ws.Range("K77").Value = 6514.8485
$ws.Range("M77").Value = -2146.8485
$ws.Range("H97").Value = 1019.1111
$ws.Range("I97").Value = 613.3333
$ws.Range("K97").Value = 613.3333
$ws.Range("M97").Value = -117.3333
$ws.Range("H102").Value = 7410567
$ws.Range("J102").Value = 5000
$ws.Range("L102").Value = 5000
$ws.Range("N102").Value = -8244
$ws.Range("H116").Value = 3582
$ws.Range("I116").Value = 3927.75
$ws.Range("J116").Value = 2660
$ws.Range("K116").Value = 3927.75
$ws.Range("L116").Value = 2660
$ws.Range("M116").Value = -1633.75
$ws.Range("N116").Value = -7248
$ws.Range("H136").Value = 5221.731
$ws.Range("I136").Value = 6036.476
$ws.Range("J136").Value = 1799.8
$ws.Range("K136").Value = 18109.428
$ws.Range("L136").Value = 5399.4
$ws.Range("M136").Value = -15559.428
$ws.Range("N136").Value = -10499.4

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H3").Value = 3582
$ws.Range("I3").Value = 3927.75
$ws.Range("J3").Value = 2660
$ws.Range("K3").Value = 3927.75
$ws.Range("L3").Value = 2660
$ws.Range("M3").Value = -3813.75
$ws.Range("N3").Value = -2888
$ws.Range("H86").Value = 15874599
$ws.Range("I86").Value = 23810984
$ws.Range("K86").Value = 23810984
$ws.Range("M86").Value = -23809861
$ws.Range("H89").Value = 15874599
$ws.Range("I89").Value = 23810984
$ws.Range("K89").Value = 119054920
$ws.Range("M89").Value = -119049304
$ws.Range("H94").Value = 2181.7
$ws.Range("I94").Value = 1816.0769
$ws.Range("K94").Value = 1816.0769
$ws.Range("M94").Value = -1365.0769
$ws.Range("H99").Value = 142858290
$ws.Range("H105").Value = 15764.2
$ws.Range("I105").Value = 26765
$ws.Range("J105").Value = 3191.8572
$ws.Range("K105").Value = 26765
$ws.Range("L105").Value = 3191.8572
$ws.Range("M105").Value = -25018
$ws.Range("N105").Value = -6685.8572
$ws.Range("H107").Value = 1034.1818
$ws.Range("I107").Value = 972.4667
$ws.Range("J107").Value = 1166.4286
$ws.Range("K107").Value = 972.4667
$ws.Range("L107").Value = 1166.4286
$ws.Range("M107").Value = 947.5333
$ws.Range("N107").Value = -5006.4286
$ws.Range("H134").Value = 4154.96
$ws.Range("I134").Value = 5337.0347
$ws.Range("J134").Value = 2522.5715
$ws.Range("K134").Value = 16011.1041
$ws.Range("L134").Value = 7567.7145
$ws.Range("M134").Value = -13476.1041
$ws.Range("N134").Value = -12637.7145

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H82").Value = 4000
$ws.Range("I82").Value = 0
$ws.Range("K82").Value = 0
$ws.Range("M82").ClearContents()
$ws.Range("H85").Value = 4000
$ws.Range("I85").Value = 0
$ws.Range("K85").Value = 0
$ws.Range("M85").ClearContents()
$ws.Range("H109").Value = 2297.3667
$ws.Range("I109").Value = 782.1
$ws.Range("J109").Value = 3055
$ws.Range("K109").Value = 2346.3
$ws.Range("L109").Value = 9165
$ws.Range("M109").Value = -1306.3
$ws.Range("N109").Value = -11245
$ws.Range("H122").Value = 955.44446
$ws.Range("I122").Value = 783.3333
$ws.Range("J122").Value = 1299.6666
$ws.Range("K122").Value = 7049.9997
$ws.Range("L122").Value = 11696.9994
$ws.Range("M122").Value = -4599.9997
$ws.Range("N122").Value = -16596.9994
$ws.Range("H125").Value = 3578.1667
$ws.Range("I125").Value = 1500
$ws.Range("J125").Value = 3837.9375
$ws.Range("K125").Value = 4500
$ws.Range("L125").Value = 11513.8125
$ws.Range("M125").Value = 420
$ws.Range("N125").Value = -21353.8125

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H80").Value = 3416.6667
$ws.Range("I80").Value = 3466.6667
$ws.Range("J80").Value = 3366.6667
$ws.Range("K80").Value = 3466.6667
$ws.Range("L80").Value = 3366.6667
$ws.Range("M80").Value = -2468.6667
$ws.Range("N80").Value = -5362.6667
$ws.Range("H83").Value = 3416.6667
$ws.Range("I83").Value = 3466.6667
$ws.Range("J83").Value = 3366.6667
$ws.Range("K83").Value = 17333.3335
$ws.Range("L83").Value = 16833.3335
$ws.Range("M83").Value = -12341.3335
$ws.Range("N83").Value = -26817.3335
$ws.Range("H102").Value = 1535
$ws.Range("I102").Value = 1396.2
$ws.Range("J102").Value = 1882
$ws.Range("K102").Value = 1396.2
$ws.Range("L102").Value = 1882
$ws.Range("M102").Value = 225.8
$ws.Range("N102").Value = -5126
$ws.Range("H113").Value = 125001200
$ws.Range("I113").Value = 333333900
$ws.Range("J113").Value = 1580
$ws.Range("K113").Value = 333333900
$ws.Range("L113").Value = 1580
$ws.Range("M113").Value = -333331730
$ws.Range("N113").Value = -5920
$ws.Range("H132").Value = 3142
$ws.Range("I132").Value = 2155
$ws.Range("J132").Value = 3765.3684
$ws.Range("K132").Value = 6465
$ws.Range("L132").Value = 11296.1052
$ws.Range("M132").Value = -3935
$ws.Range("N132").Value = -16356.1052
$ws.Range("H134").Value = 34326
$ws.Range("J134").Value = 34326
$ws.Range("L134").Value = 102978
$ws.Range("N134").Value = -108048
